# Apply cryptos price/volume update (generated from OOXML diff)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '62.567.79'
$ws.Range("E2").Value = '  -0.94%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.454.43'
$ws.Range("E3").Value = '  -0.32%  '
$ws.Range("E4").Value = '  +0.10%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '570.03'
$ws.Range("E5").Value = '  -1.34%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '146.74'
$ws.Range("E6").Value = '  +0.10%  '
$ws.Range("E7").Value = '  +0.02%  '
$ws.Range("E8").Value = '  -1.82%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.111'
$ws.Range("E9").Value = '  -0.51%  '
$ws.Range("E10").Value = '  +0.05%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '5.26'
$ws.Range("E11").Value = '  -0.44%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.351'
$ws.Range("E12").Value = '  -1.24%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '28.74'
$ws.Range("E13").Value = '  -1.27%  '
$ws.Range("E14").Value = '  -2.86%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '2.901.40'
$ws.Range("E15").Value = '  -0.27%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '62.527.77'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.456.34'
$ws.Range("E17").Value = '  -0.05%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '7.91'
$ws.Range("E18").Value = '  -0.02%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '10.83'
$ws.Range("E19").Value = '  -2.48%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '324.40'
$ws.Range("E20").Value = '  -1.91%  '
$ws.Range("E21").Value = '  -0.11%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '2.17'
$ws.Range("E22").Value = '  -3.03%  '
$ws.Range("E23").Value = '  -0.11%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '10.00'
$ws.Range("E24").Value = '  +9.54%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '65.30'
$ws.Range("E25").Value = '  -1.92%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '638.18'
$ws.Range("E26").Value = '  -4.35%  '
$ws.Range("E27").Value = '  +0.04%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.0₃0964'
$ws.Range("E28").Value = '  -4.49%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.00'
$ws.Range("E29").Value = '  -5.54%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.42'
$ws.Range("E30").Value = '  -2.27%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.90'
$ws.Range("E31").Value = '  -3.35%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.82'
$ws.Range("E32").Value = '  -3.71%  '
$ws.Range("E33").Value = '  -4.80%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.999'
$ws.Range("E34").Value = '  -0.01%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.50'
$ws.Range("E35").Value = '  -3.69%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '4.72'
$ws.Range("E36").Value = '  -1.73%  '
$ws.Range("B37").Value = 'Monero'
$ws.Range("C37").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '150.94'
$ws.Range("E37").Value = '  -1.20%  '
$ws.Range("B38").Value = 'PolygonEcosystemToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.366'
$ws.Range("E38").Value = '  -2.19%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '18.54'
$ws.Range("E39").Value = '  -1.75%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '5.27'
$ws.Range("E40").Value = '  -5.64%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.69'
$ws.Range("E41").Value = '  -2.02%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.72'
$ws.Range("E42").Value = '  -2.74%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0₆0309'
$ws.Range("E43").Value = '  +1.93%  '
$ws.Range("E44").Value = '  -0.13%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '152.03'
$ws.Range("E45").Value = '  +3.77%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '15.29'
$ws.Range("E46").Value = '  +1.06%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.56'
$ws.Range("E47").Value = '  -2.06%  '
$ws.Range("E48").Value = '  -0.41%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '20.18'
$ws.Range("E49").Value = '  -3.00%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0505'
$ws.Range("E50").Value = '  -2.48%  '
$ws.Range("E51").Value = '  -1.78%  '
